# fix(bar chart): remove jspdf
# Clears the sample row of generated data (uuid / instance type / region /
# metrics / pricing model) from Sheet1 row 2, leaving only the formatting
# (styles) on A2 and J2 intact, and moves the active selection to E10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the sample data values that were previously filled in on row 2,
# while keeping the existing cell formatting (styles) untouched.
$ws.Range("A2").ClearContents()
$ws.Range("C2:J2").ClearContents()

# Update the saved cursor/selection position on Sheet1.
$ws.Range("E10").Select()
